# guion tema 1 grado 08
# Update the "RECURSOS APROVECHADOS" sheet: the MATERIA changed from MA to MT,
# and the GUION (recurso) title changed accordingly; also move the active
# selection and drop the stale "highlighted" formatting on the old D2 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RECURSOS APROVECHADOS")

# MATERIA column (B): MA -> MT for every data row
$ws.Range("B2:B10").Value2 = "MT"

# GUION column (C): rows 2-8 belong to the "racionales e irracionales" guion,
# rows 9-10 belong to the "numeros reales" guion
$ws.Range("C2:C8").Value2 = "Los números racionales e irracionales"
$ws.Range("C9:C10").Value2 = "Los números reales"

# Row 2 used to carry a one-off highlighted look (custom row height + a
# Times New Roman style on D2); both get cleared back to the sheet defaults.
$ws.Range("D2").Style = "Normal"
$ws.Rows.Item(2).AutoFit()

# Move the saved selection from C24 to C16
$ws.Activate()
$ws.Range("C16").Select()
